$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows starting at row 404 (pushes old rows 404-424 down to 408-428)
$ws.Rows.Item(404).Resize(4).Insert()

# Shared (unchanged) column values for all the new "Sandia" rows
$A = 9
$B = "Vega Central Mapocho de Santiago"
$C = "Metropolitana"
$E = 13
$F = 100112028
$G = "Sandia"
$H = "Sin especificar"
$N = '$/unidad'
$O = "Región de O'Higgins"
$Q = 1
$R = "Hortaliza"

$newRows = @(
    @{ Row = 404; D = 44610; I = "Extra";   J = 250; K = 2500; L = 2800; M = 2650; P = 2650 },
    @{ Row = 405; D = 44610; I = "Primera"; J = 520; K = 2000; L = 2300; M = 2150; P = 2150 },
    @{ Row = 406; D = 44610; I = "Segunda"; J = 340; K = 1600; L = 1900; M = 1750; P = 1750 },
    @{ Row = 407; D = 44610; I = "Tercera"; J = 160; K = 1400; L = 1400; M = 1400; P = 1400 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $A
    $ws.Cells.Item($row, 2).Value = $B
    $ws.Cells.Item($row, 3).Value = $C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $E
    $ws.Cells.Item($row, 6).Value = $F
    $ws.Cells.Item($row, 7).Value = $G
    $ws.Cells.Item($row, 8).Value = $H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $N
    $ws.Cells.Item($row, 15).Value = $O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $Q
    $ws.Cells.Item($row, 18).Value = $R
}
